$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 41
$ws1.Range("F4").Value = 3495
$ws1.Range("G4").Value = 67.5
$ws1.Range("F5").Value = 3495
$ws1.Range("G5").Value = 67.5
$ws1.Range("F6").Value = 250
$ws1.Range("F7").Value = 5027
$ws1.Range("F11").Value = 670
$ws1.Range("F12").Value = 291
$ws1.Range("F13").Value = 72
$ws1.Range("F14").Value = 29
$ws1.Range("F16").Value = 304
$ws1.Range("F22").Value = 4867
$ws1.Range("F26").Value = 5985
$ws1.Range("F29").Value = 3213
$ws1.Range("F30").Value = 323
$ws1.Range("F32").Value = 4439
$ws1.Range("F34").Value = 117
$ws1.Range("F36").Value = 968
$ws1.Range("F40").Value = 849
$ws1.Range("F41").Value = 941

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 47

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 41
$ws4.Range("F5").Value = 47
$ws4.Range("F8").Value = 3495
$ws4.Range("G8").Value = 67.5
$ws4.Range("F9").Value = 3495
$ws4.Range("G9").Value = 67.5
$ws4.Range("F10").Value = 250
$ws4.Range("F11").Value = 5027
$ws4.Range("F15").Value = 670
$ws4.Range("F16").Value = 72
$ws4.Range("F17").Value = 29
$ws4.Range("F19").Value = 304
$ws4.Range("F26").Value = 4867
$ws4.Range("F30").Value = 5985
$ws4.Range("F33").Value = 3213
$ws4.Range("F34").Value = 323
$ws4.Range("F36").Value = 4439
$ws4.Range("F39").Value = 117
$ws4.Range("F41").Value = 968
$ws4.Range("F45").Value = 849
$ws4.Range("F46").Value = 941
